# Generate Report for Handoff
# Flip the per-language handback status back to "Ready for handoff" and
# stamp the refreshed handoff timestamps, then tighten the now-shorter
# status columns back down to fit.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus          # zh-cn status
$overview.Range("F2").Value = $newStatus          # de-de status
$overview.Range("G2").Value = "2016-10-25 03:03:22"  # Latest HO Xliff Generate Date

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus               # Status
$zhcn.Range("H2").Value = "2016-10-25 03:03:09"    # Latest Handoff Datetime

# --- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus               # Status
$dede.Range("H2").Value = "2016-10-25 03:03:22"    # Latest Handoff Datetime

# --- Re-fit the status columns now that the text is shorter ----------
# (ColumnWidth is in characters; Excel quantizes it to whole pixels, so
# 16.33 is the input that lands closest to the ~17.22-char target width.)
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
